$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (date + 24 hourly prices + daily average)
$ws.Range("A2").Value = 45939
$ws.Range("B2").Value = 92.37
$ws.Range("C2").Value = 80.95999999999999
$ws.Range("D2").Value = 78.95
$ws.Range("E2").Value = 69.86
$ws.Range("F2").Value = 65.58
$ws.Range("G2").Value = 68.47
$ws.Range("H2").Value = 85.78
$ws.Range("I2").Value = 101.88
$ws.Range("J2").Value = 111.28
$ws.Range("K2").Value = 91.06
$ws.Range("L2").Value = 58.03
$ws.Range("M2").Value = 23.33
$ws.Range("N2").Value = 19.34
$ws.Range("O2").Value = 19.01
$ws.Range("P2").Value = 19.01
$ws.Range("Q2").Value = 22.03
$ws.Range("R2").Value = 39.55
$ws.Range("S2").Value = 55.01
$ws.Range("T2").Value = 75.51000000000001
$ws.Range("U2").Value = 106
$ws.Range("V2").Value = 114.46
$ws.Range("W2").Value = 101.69
$ws.Range("X2").Value = 95.06999999999999
$ws.Range("Y2").Value = 92.73999999999999
$ws.Range("Z2").Value = 70.29000000000001

# Slot summary columns
$ws.Range("AB2").Value = 100.99
$ws.Range("AD2").Value = 108.07
$ws.Range("AF2").Value = 101.17
$ws.Range("AG2").Value = "3h-17h"
